$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
# (the source workbook stores these columns as text, e.g. "42.713.18",
# "2.580.39", "1.00" etc.) -- force text format first so Excel does not
# auto-convert the literal into a floating point number.
$textForceCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D13", "D16", "D17", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D31", "D32", "D33", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D50")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume/change columns, and the
# Cosmos/Dai row swap in rows 26-27).
$ws.Range("D2").Value = "42.721.53"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "2.251.09"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "248.97"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "70.05"
$ws.Range("E7").Value = "  +6.18%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +17.83%  "
$ws.Range("D10").Value = "39.18"
$ws.Range("E10").Value = "  +10.43%  "
$ws.Range("D11").Value = "59.35"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "0.0964"
$ws.Range("E12").Value = "  +4.35%  "
$ws.Range("D13").Value = "7.49"
$ws.Range("E13").Value = "  +8.47%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "2.580.53"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("D16").Value = "14.81"
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("D17").Value = "0.880"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("D18").Value = "2.236.65"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "42.647.33"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").Value = "6.29"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("D22").Value = "72.86"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").Value = "235.24"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("E25").Value = "  +6.06%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "11.39"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "167.36"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "20.91"
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("D33").Value = "6.43"
$ws.Range("E33").Value = "  +14.91%  "
$ws.Range("E34").Value = "  +6.46%  "
$ws.Range("D35").Value = "0.0793"
$ws.Range("E35").Value = "  +7.86%  "
$ws.Range("D36").Value = "31.09"
$ws.Range("E36").Value = "  +22.68%  "
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("D38").Value = "4.41"
$ws.Range("E38").Value = "  +9.96%  "
$ws.Range("D39").Value = "4.69"
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("D40").Value = "0.0320"
$ws.Range("E40").Value = "  +7.18%  "
$ws.Range("D41").Value = "2.31"
$ws.Range("E41").Value = "  +6.41%  "
$ws.Range("D42").Value = "12.53"
$ws.Range("E42").Value = "  +7.99%  "
$ws.Range("E43").Value = "  +5.98%  "
$ws.Range("D44").Value = "62.27"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("E45").Value = "  +6.93%  "
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").Value = "  +6.14%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("E51").Value = "  +4.09%  "

# Restore default styling on the cells we temporarily reformatted as text
# so only the values (not the number format) differ from the original.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
